# Add a new worksheet "ViewUserPage" right after "UserPage" and before
# "ResetPage", populate it with the "view user" detail labels, and make it
# the active/selected sheet (mirrors the author's commit "9 testcases").

$wb = $excel.ActiveWorkbook

$userPage = $wb.Worksheets.Item("UserPage")

$viewUserPage = $wb.Worksheets.Add($null, $userPage)
$viewUserPage.Name = "ViewUserPage"

# Column A width -> stored OOXML width of 34 (Excel pads ColumnWidth by
# ~0.8333 when serializing, so compensate for it here).
$viewUserPage.Columns.Item(1).ColumnWidth = 33.16666666666667

$viewUserPage.Range("A1").Value = "Email:"
$viewUserPage.Range("A2").Value = "Role:"
$viewUserPage.Range("A3").Value = "Username:"
$viewUserPage.Range("A4").Value = "Sales Commission Percentage (%):"
$viewUserPage.Range("A5").Value = "Active"
$viewUserPage.Range("A6").Value = "Sales Commission Percentage (%):"
$viewUserPage.Range("A7").Value = "Allowed Contacts: All"

$viewUserPage.Range("A7").Select()
